# Apply cryptocurrency price/volume update (GitHub Actions refresh)
# Rows 32/33 and 37/38 also swap B (Coin) / C (Link) since the source
# ranking reordered those coin pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A literal apostrophe as a one-character string (kept in a variable to avoid
# quoting headaches inside the double-quoted literals built below).
$apos = "'"

# Value set per row: B (coin name), C (link), D (price text), E (volume text)
$updates = @(
    @{ Row = 2; B = $null; C = $null; D = "19.994.50"; E = "  -8.10%  " }
    @{ Row = 3; B = $null; C = $null; D = "1.419.42"; E = "  -7.82%  " }
    @{ Row = 4; B = $null; C = $null; D = "1.001"; E = "  +0.02%  " }
    @{ Row = 5; B = $null; C = $null; D = "1.001"; E = "  +0.04%  " }
    @{ Row = 6; B = $null; C = $null; D = "273.52"; E = "  -5.74%  " }
    @{ Row = 7; B = $null; C = $null; D = "0.3728"; E = "  -4.03%  " }
    @{ Row = 8; B = $null; C = $null; D = "0.3071"; E = "  -3.74%  " }
    @{ Row = 9; B = $null; C = $null; D = "39.69"; E = "  -8.09%  " }
    @{ Row = 10; B = $null; C = $null; D = $null; E = "  -4.54%  " }
    @{ Row = 11; B = $null; C = $null; D = "0.06602"; E = "  -8.37%  " }
    @{ Row = 12; B = $null; C = $null; D = "1.001"; E = "  +0.02%  " }
    @{ Row = 13; B = $null; C = $null; D = "5.417"; E = "  -3.94%  " }
    @{ Row = 14; B = $null; C = $null; D = "17.09"; E = "  -8.12%  " }
    @{ Row = 15; B = $null; C = $null; D = "6.162"; E = "  -6.67%  " }
    @{ Row = 16; B = $null; C = $null; D = "1.421.96"; E = "  -7.98%  " }
    @{ Row = 17; B = $null; C = $null; D = "0.00001008"; E = "  -9.41%  " }
    @{ Row = 18; B = $null; C = $null; D = "0.05826"; E = "  -11.53%  " }
    @{ Row = 19; B = $null; C = $null; D = "74.63"; E = "  -10.50%  " }
    @{ Row = 20; B = $null; C = $null; D = "1.001"; E = "  +0.11%  " }
    @{ Row = 21; B = $null; C = $null; D = "5.642"; E = "  -8.19%  " }
    @{ Row = 22; B = $null; C = $null; D = "14.45"; E = "  -6.14%  " }
    @{ Row = 23; B = $null; C = $null; D = "10.95"; E = "  -0.29%  " }
    @{ Row = 24; B = $null; C = $null; D = "2.324"; E = "  -2.87%  " }
    @{ Row = 25; B = $null; C = $null; D = "20.002.91"; E = "  -8.09%  " }
    @{ Row = 26; B = $null; C = $null; D = "2.287"; E = "  -4.10%  " }
    @{ Row = 27; B = $null; C = $null; D = "138.75"; E = "  -5.41%  " }
    @{ Row = 28; B = $null; C = $null; D = "16.92"; E = "  -7.96%  " }
    @{ Row = 29; B = $null; C = $null; D = "1.582.00"; E = "  -7.95%  " }
    @{ Row = 30; B = $null; C = $null; D = "108.98"; E = "  -7.24%  " }
    @{ Row = 31; B = $null; C = $null; D = "3.811"; E = "  -21.38%  " }
    @{ Row = 32; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "0.8872"; E = "  -8.66%  " }
    @{ Row = 33; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "5.419"; E = "  -8.26%  " }
    @{ Row = 34; B = $null; C = $null; D = "0.07733"; E = "  -5.78%  " }
    @{ Row = 35; B = $null; C = $null; D = "8.449"; E = "  -5.19%  " }
    @{ Row = 36; B = $null; C = $null; D = $null; E = "  +5.34%  " }
    @{ Row = 37; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "4.781"; E = "  -7.25%  " }
    @{ Row = 38; B = "Frax"; C = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D = "1.001"; E = "  +0.10%  " }
    @{ Row = 39; B = $null; C = $null; D = "0.05683"; E = "  -6.53%  " }
    @{ Row = 40; B = $null; C = $null; D = "0.1915"; E = "  -6.16%  " }
    @{ Row = 41; B = $null; C = $null; D = "0.02029"; E = "  -8.03%  " }
    @{ Row = 42; B = $null; C = $null; D = "1.083"; E = "  -9.12%  " }
    @{ Row = 43; B = $null; C = $null; D = "1.261"; E = "  -15.21%  " }
    @{ Row = 44; B = $null; C = $null; D = "0.5317"; E = "  -7.68%  " }
    @{ Row = 45; B = $null; C = $null; D = "3.533"; E = "  -5.72%  " }
    @{ Row = 46; B = $null; C = $null; D = "12.26"; E = "  -5.60%  " }
    @{ Row = 47; B = $null; C = $null; D = "0.5130"; E = "  -7.18%  " }
    @{ Row = 48; B = $null; C = $null; D = "1.800"; E = "  -3.81%  " }
    @{ Row = 49; B = $null; C = $null; D = "109.54"; E = "  -7.29%  " }
    @{ Row = 50; B = $null; C = $null; D = "1.048"; E = "  -8.49%  " }
    @{ Row = 51; B = $null; C = $null; D = "1.001"; E = "  +0.04%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.B) { $ws.Range("B$row").Value = $u.B }
    if ($null -ne $u.C) { $ws.Range("C$row").Value = $u.C }
    if ($null -ne $u.D) {
        # Price column holds text, not numbers (e.g. "19.994.50"). Force text
        # with a leading apostrophe so Excel does not coerce look-alike values
        # such as "1.001" or "0.8872" into numeric cells.
        $ws.Range("D$row").Value = $apos + $u.D
    }
    if ($null -ne $u.E) { $ws.Range("E$row").Value = $u.E }
}
